# Fix: Production sync issues, Overtime logic, and Type Safety Audit improvements
$wb = $excel.ActiveWorkbook

$wsWorkforce = $wb.Worksheets.Item("WORKFORCE PLANNING")
$wsComp      = $wb.Worksheets.Item("COMPENSATION STRATEGY")
$wsLabor     = $wb.Worksheets.Item("LABOR COST ANALYSIS")
$wsCross     = $wb.Worksheets.Item("CROSS REFERENCE")

# --- WORKFORCE PLANNING ---------------------------------------------------
# Benchmark note text should reflect the new hiring-fee benchmark.
$wsWorkforce.Range("C5").Value = 'Benchmark (Sales Hire): $240'

# Est. Severance (per worker): 2000 -> 220
$wsWorkforce.Range("B6").Value = 220

# Overtime cost constant used in the HIRE vs OVERTIME optimizer: 650 -> 27.3
$wsWorkforce.Range("B21").Formula = '=B19*($B$5 + 27.3)'
$wsWorkforce.Range("B24").Formula = '=MIN(B19, B23) * 27.3 * 1.4'

# --- COMPENSATION STRATEGY -------------------------------------------------
# Zone base salary + proposed salary values, and related strike-risk formula
# 27.3 replaces 650 as the wage floor; 28 replaces 676 as the proposed salary.
$zoneRows = 11..15
foreach ($row in $zoneRows) {
    $wsComp.Range("B$row").Value = 27.3
    $wsComp.Range("D$row").Value = 28
    $wsComp.Range("E$row").Formula = "=IF(OR(D$row<C$row, D$row<27.3),""STRIKE RISK!"",""OK"")"
}

# --- LABOR COST ANALYSIS ---------------------------------------------------
# Repair cross-sheet references to use the quoted real sheet name
# (sheet tab is "WORKFORCE PLANNING", with a space).
$wsLabor.Range("B9").Formula = "='WORKFORCE PLANNING'!C15"
$wsLabor.Range("B14").Formula = "='WORKFORCE PLANNING'!K15"

# --- CROSS REFERENCE --------------------------------------------------------
# Leading apostrophe forces this numeric-looking literal to stay text (it was
# already stored as text before the edit).
$wsCross.Range("B5").Value = "'0.0"
$wsCross.Range("B11").Value = "🟢 OK"
